$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# and mangles formatting (e.g. "1.00" -> 1, "304.96" -> 304.9599999...).
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D17", "D19", "D21", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D37", "D40", "D43", "D44", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.087.39"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "2.305.21"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "304.96"
$ws.Range("E5").Value = "  +1.77%  "

$ws.Range("D6").Value = "97.22"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  -1.16%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").Value = "35.52"
$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").Value = "18.62"
$ws.Range("E12").Value = "  +5.32%  "

$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").Value = "6.91"
$ws.Range("E14").Value = "  +2.47%  "

$ws.Range("D15").Value = "2.664.16"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").Value = "2.317.94"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").Value = "43.008.02"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").Value = "12.61"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("D23").Value = "236.84"
$ws.Range("E23").Value = "  -1.74%  "

$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").Value = "25.01"
$ws.Range("E27").Value = "  -0.46%  "

$ws.Range("D28").Value = "2.18"
$ws.Range("E28").Value = "  +7.56%  "

$ws.Range("D29").Value = "166.02"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").Value = "9.05"
$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("D31").Value = "33.00"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("D33").Value = "18.19"
$ws.Range("E33").Value = "  +6.65%  "

$ws.Range("D34").Value = "5.00"
$ws.Range("E34").Value = "  -0.17%  "

$ws.Range("D35").Value = "4.46"
$ws.Range("E35").Value = "  -7.38%  "

$ws.Range("E36").Value = "  -0.91%  "

$ws.Range("D37").Value = "0.0689"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").Value = "1.997.49"
$ws.Range("E42").Value = "  -0.49%  "

$ws.Range("D43").Value = "10.46"
$ws.Range("E43").Value = "  +3.46%  "

$ws.Range("D44").Value = "0.0280"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.11"
$ws.Range("E45").Value = "  +2.46%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "17.91"
$ws.Range("E46").Value = "  +4.93%  "

$ws.Range("D47").Value = "2.79"
$ws.Range("E47").Value = "  +0.60%  "

$ws.Range("D48").Value = "53.64"
$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("D49").Value = "2.532.06"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("E50").Value = "  -7.46%  "

$ws.Range("D51").Value = "71.73"
$ws.Range("E51").Value = "  -0.26%  "
